$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: the old "k@gmail.com" shared string is being replaced in-place by
# "XYZ123" (so B3, which used to reference the old "XYZ123" string slot,
# now points at the same "XYZ123" text while A3 becomes the new
# "prashant" display text for the existing rId2 hyperlink, now showing
# "k@gmail.com" as the stored display attribute).
$ws.Range("A3").Value = "prashant"
$ws.Range("B3").Value = "XYZ123"

# Row 4: a brand-new step/row for the "not registered" email, with its own
# mailto hyperlink.
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:prashant@test.com", "", "", "prashant@test.com")
$ws.Range("A4").Style = "Hyperlink"
